$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sELE")
$ws.Cells.Clear()
$ws.Cells.Item(1, 1).Value = 'Sector'
$ws.Cells.Item(1, 2).Value = 'Year'
$ws.Cells.Item(1, 3).Value = 'COL'
$ws.Cells.Item(1, 4).Value = 'OIL'
$ws.Cells.Item(1, 5).Value = 'GAS'
$ws.Cells.Item(1, 6).Value = 'BMS'
$ws.Cells.Item(1, 7).Value = 'ELE'
$ws.Cells.Item(2, 1).Value = 'Industry'
$ws.Cells.Item(2, 2).Value = 2010
$ws.Cells.Item(2, 3).Value = 0.235317279731797
$ws.Cells.Item(2, 4).Value = 0.2542029283506
$ws.Cells.Item(2, 5).Value = 0.118208827636883
$ws.Cells.Item(2, 6).Value = 0.0364559087844255
$ws.Cells.Item(2, 7).Value = 0.355815055496294
$ws.Cells.Item(3, 1).Value = 'Industry'
$ws.Cells.Item(3, 2).Value = 2011
$ws.Cells.Item(3, 3).Value = 0.233781164697422
$ws.Cells.Item(3, 4).Value = 0.265632827450878
$ws.Cells.Item(3, 5).Value = 0.123996041168318
$ws.Cells.Item(3, 6).Value = 0.0393771825968263
$ws.Cells.Item(3, 7).Value = 0.337212784086556
$ws.Cells.Item(4, 1).Value = 'Industry'
$ws.Cells.Item(4, 2).Value = 2012
$ws.Cells.Item(4, 3).Value = 0.238183139592608
$ws.Cells.Item(4, 4).Value = 0.254729802181916
$ws.Cells.Item(4, 5).Value = 0.126204574849316
$ws.Cells.Item(4, 6).Value = 0.0385824030204225
$ws.Cells.Item(4, 7).Value = 0.342300080355736
$ws.Cells.Item(5, 1).Value = 'Industry'
$ws.Cells.Item(5, 2).Value = 2013
$ws.Cells.Item(5, 3).Value = 0.238991514104705
$ws.Cells.Item(5, 4).Value = 0.25628028219994
$ws.Cells.Item(5, 5).Value = 0.118880058675235
$ws.Cells.Item(5, 6).Value = 0.0418285950752597
$ws.Cells.Item(5, 7).Value = 0.34401954994486
$ws.Cells.Item(6, 1).Value = 'Industry'
$ws.Cells.Item(6, 2).Value = 2014
$ws.Cells.Item(6, 3).Value = 0.239516304821693
$ws.Cells.Item(6, 4).Value = 0.252808385648269
$ws.Cells.Item(6, 5).Value = 0.119951336210655
$ws.Cells.Item(6, 6).Value = 0.0394015934257564
$ws.Cells.Item(6, 7).Value = 0.348322379893626
$ws.Cells.Item(7, 1).Value = 'Industry'
$ws.Cells.Item(7, 2).Value = 2015
$ws.Cells.Item(7, 3).Value = 0.244278082479678
$ws.Cells.Item(7, 4).Value = 0.240704426350651
$ws.Cells.Item(7, 5).Value = 0.123588729268783
$ws.Cells.Item(7, 6).Value = 0.040869790218956
$ws.Cells.Item(7, 7).Value = 0.350558971681932
$ws.Cells.Item(8, 1).Value = 'Industry'
$ws.Cells.Item(8, 2).Value = 2016
$ws.Cells.Item(8, 3).Value = 0.242373154334521
$ws.Cells.Item(8, 4).Value = 0.228317868362445
$ws.Cells.Item(8, 5).Value = 0.130553448532646
$ws.Cells.Item(8, 6).Value = 0.0424868181183189
$ws.Cells.Item(8, 7).Value = 0.356268710652069
$ws.Cells.Item(9, 1).Value = 'Industry'
$ws.Cells.Item(9, 2).Value = 2017
$ws.Cells.Item(9, 3).Value = 0.24922368028528
$ws.Cells.Item(9, 4).Value = 0.217978539482308
$ws.Cells.Item(9, 5).Value = 0.130354331207477
$ws.Cells.Item(9, 6).Value = 0.0435277511261819
$ws.Cells.Item(9, 7).Value = 0.358915697898754
$ws.Cells.Item(10, 1).Value = 'Industry'
$ws.Cells.Item(10, 2).Value = 2018
$ws.Cells.Item(10, 3).Value = 0.242436419524306
$ws.Cells.Item(10, 4).Value = 0.217749202182963
$ws.Cells.Item(10, 5).Value = 0.136550016297178
$ws.Cells.Item(10, 6).Value = 0.0452056222686474
$ws.Cells.Item(10, 7).Value = 0.358058739726906
$ws.Cells.Item(11, 1).Value = 'Industry'
$ws.Cells.Item(11, 2).Value = 2019
$ws.Cells.Item(11, 3).Value = 0.246172595109221
$ws.Cells.Item(11, 4).Value = 0.21124427442343
$ws.Cells.Item(11, 5).Value = 0.134886409911155
$ws.Cells.Item(11, 6).Value = 0.0459790921710079
$ws.Cells.Item(11, 7).Value = 0.361717628385187
$ws.Cells.Item(12, 1).Value = 'Industry'
$ws.Cells.Item(12, 2).Value = 2020
$ws.Cells.Item(12, 3).Value = 0.2378380422764
$ws.Cells.Item(12, 4).Value = 0.218958939914931
$ws.Cells.Item(12, 5).Value = 0.134617175435078
$ws.Cells.Item(12, 6).Value = 0.0445472571891143
$ws.Cells.Item(12, 7).Value = 0.364038585184477
$ws.Cells.Item(13, 1).Value = 'Transport'
$ws.Cells.Item(13, 2).Value = 2010
$ws.Cells.Item(13, 3).Value = 0.0000127563610961359
$ws.Cells.Item(13, 4).Value = 0.976242188627097
$ws.Cells.Item(13, 5).Value = 0.00127593983249683
$ws.Cells.Item(13, 6).Value = 0.00252089993090304
$ws.Cells.Item(13, 7).Value = 0.0199482152484074
$ws.Cells.Item(14, 1).Value = 'Transport'
$ws.Cells.Item(14, 2).Value = 2011
$ws.Cells.Item(14, 3).Value = 0.0000130761890381062
$ws.Cells.Item(14, 4).Value = 0.97628601984156
$ws.Cells.Item(14, 5).Value = 0.00128582525541378
$ws.Cells.Item(14, 6).Value = 0.00251280766015607
$ws.Cells.Item(14, 7).Value = 0.0199022710538319
$ws.Cells.Item(15, 1).Value = 'Transport'
$ws.Cells.Item(15, 2).Value = 2012
$ws.Cells.Item(15, 3).Value = 0.0000122700235773222
$ws.Cells.Item(15, 4).Value = 0.976158085725306
$ws.Cells.Item(15, 5).Value = 0.00122448542981893
$ws.Cells.Item(15, 6).Value = 0.00258677266289086
$ws.Cells.Item(15, 7).Value = 0.0200183861584066
$ws.Cells.Item(16, 1).Value = 'Transport'
$ws.Cells.Item(16, 2).Value = 2013
$ws.Cells.Item(16, 3).Value = 0.0000120951891385202
$ws.Cells.Item(16, 4).Value = 0.97548941750524
$ws.Cells.Item(16, 5).Value = 0.00113121847890265
$ws.Cells.Item(16, 6).Value = 0.00306135602984965
$ws.Cells.Item(16, 7).Value = 0.0203059127968693
$ws.Cells.Item(17, 1).Value = 'Transport'
$ws.Cells.Item(17, 2).Value = 2014
$ws.Cells.Item(17, 3).Value = 0.0000124358001815627
$ws.Cells.Item(17, 4).Value = 0.974323345204068
$ws.Cells.Item(17, 5).Value = 0.00109827751077169
$ws.Cells.Item(17, 6).Value = 0.00382662661902664
$ws.Cells.Item(17, 7).Value = 0.0207393148659519
$ws.Cells.Item(18, 1).Value = 'Transport'
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 0.0000121347680783447
$ws.Cells.Item(18, 4).Value = 0.97356260392445
$ws.Cells.Item(18, 5).Value = 0.000978980613888078
$ws.Cells.Item(18, 6).Value = 0.00472436038293392
$ws.Cells.Item(18, 7).Value = 0.0207219203106501
$ws.Cells.Item(19, 1).Value = 'Transport'
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 0.0000132625027271021
$ws.Cells.Item(19, 4).Value = 0.972840052227736
$ws.Cells.Item(19, 5).Value = 0.000840511110330097
$ws.Cells.Item(19, 6).Value = 0.00536667172852187
$ws.Cells.Item(19, 7).Value = 0.0209395024306852
$ws.Cells.Item(20, 1).Value = 'Transport'
$ws.Cells.Item(20, 2).Value = 2017
$ws.Cells.Item(20, 3).Value = 0.0000133731630706788
$ws.Cells.Item(20, 4).Value = 0.972273086676485
$ws.Cells.Item(20, 5).Value = 0.000719476173202521
$ws.Cells.Item(20, 6).Value = 0.00587984547310072
$ws.Cells.Item(20, 7).Value = 0.0211142185141413
$ws.Cells.Item(21, 1).Value = 'Transport'
$ws.Cells.Item(21, 2).Value = 2018
$ws.Cells.Item(21, 3).Value = 0.0000121864816036597
$ws.Cells.Item(21, 4).Value = 0.972229378018905
$ws.Cells.Item(21, 5).Value = 0.000593752464800533
$ws.Cells.Item(21, 6).Value = 0.00607733067306953
$ws.Cells.Item(21, 7).Value = 0.0210873523616217
$ws.Cells.Item(22, 1).Value = 'Transport'
$ws.Cells.Item(22, 2).Value = 2019
$ws.Cells.Item(22, 3).Value = 0.0000124402436352604
$ws.Cells.Item(22, 4).Value = 0.972049191487832
$ws.Cells.Item(22, 5).Value = 0.000483441690159145
$ws.Cells.Item(22, 6).Value = 0.00592846721684796
$ws.Cells.Item(22, 7).Value = 0.0215264593615252
$ws.Cells.Item(23, 1).Value = 'Transport'
$ws.Cells.Item(23, 2).Value = 2020
$ws.Cells.Item(23, 3).Value = 0.0000137654722952663
$ws.Cells.Item(23, 4).Value = 0.968745495153598
$ws.Cells.Item(23, 5).Value = 0.000406463806940781
$ws.Cells.Item(23, 6).Value = 0.00701465525712946
$ws.Cells.Item(23, 7).Value = 0.0238196203100367
$ws.Cells.Item(24, 1).Value = 'Residential'
$ws.Cells.Item(24, 2).Value = 2010
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0.295577388456578
$ws.Cells.Item(24, 5).Value = 0.185649697865496
$ws.Cells.Item(24, 6).Value = 0.000194660599825144
$ws.Cells.Item(24, 7).Value = 0.518578253078101
$ws.Cells.Item(25, 1).Value = 'Residential'
$ws.Cells.Item(25, 2).Value = 2011
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0.294781866130855
$ws.Cells.Item(25, 5).Value = 0.192891623759889
$ws.Cells.Item(25, 6).Value = 0.000179329398467914
$ws.Cells.Item(25, 7).Value = 0.512147180710787
$ws.Cells.Item(26, 1).Value = 'Residential'
$ws.Cells.Item(26, 2).Value = 2012
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0.292714418499285
$ws.Cells.Item(26, 5).Value = 0.190718860799803
$ws.Cells.Item(26, 6).Value = 0.000168908764924207
$ws.Cells.Item(26, 7).Value = 0.516397811935988
$ws.Cells.Item(27, 1).Value = 'Residential'
$ws.Cells.Item(27, 2).Value = 2013
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0.288747412097999
$ws.Cells.Item(27, 5).Value = 0.191014059935836
$ws.Cells.Item(27, 6).Value = 0.000173333992682247
$ws.Cells.Item(27, 7).Value = 0.520065193973483
$ws.Cells.Item(28, 1).Value = 'Residential'
$ws.Cells.Item(28, 2).Value = 2014
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0.282634731811053
$ws.Cells.Item(28, 5).Value = 0.199357563841297
$ws.Cells.Item(28, 6).Value = 0.000180578946725493
$ws.Cells.Item(28, 7).Value = 0.517827125400925
$ws.Cells.Item(29, 1).Value = 'Residential'
$ws.Cells.Item(29, 2).Value = 2015
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 0.276882966712301
$ws.Cells.Item(29, 5).Value = 0.196515236635559
$ws.Cells.Item(29, 6).Value = 0.000185446790406947
$ws.Cells.Item(29, 7).Value = 0.526416349861733
$ws.Cells.Item(30, 1).Value = 'Residential'
$ws.Cells.Item(30, 2).Value = 2016
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0.275029440633412
$ws.Cells.Item(30, 5).Value = 0.200188136813701
$ws.Cells.Item(30, 6).Value = 0.00018519547382262
$ws.Cells.Item(30, 7).Value = 0.524597227079064
$ws.Cells.Item(31, 1).Value = 'Residential'
$ws.Cells.Item(31, 2).Value = 2017
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0.282907320193399
$ws.Cells.Item(31, 5).Value = 0.201477249685925
$ws.Cells.Item(31, 6).Value = 0.000177679670728215
$ws.Cells.Item(31, 7).Value = 0.515437750449948
$ws.Cells.Item(32, 1).Value = 'Residential'
$ws.Cells.Item(32, 2).Value = 2018
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 0.261349883050479
$ws.Cells.Item(32, 5).Value = 0.204412917935886
$ws.Cells.Item(32, 6).Value = 0.000192649405279954
$ws.Cells.Item(32, 7).Value = 0.534044549608355
$ws.Cells.Item(33, 1).Value = 'Residential'
$ws.Cells.Item(33, 2).Value = 2019
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 0.272207345338359
$ws.Cells.Item(33, 5).Value = 0.208935466078323
$ws.Cells.Item(33, 6).Value = 0.000194336241708559
$ws.Cells.Item(33, 7).Value = 0.518662852341609
$ws.Cells.Item(34, 1).Value = 'Residential'
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 0.267849508803155
$ws.Cells.Item(34, 5).Value = 0.212616037620618
$ws.Cells.Item(34, 6).Value = 0.00018539845399501
$ws.Cells.Item(34, 7).Value = 0.519349055122232
$ws.Cells.Item(35, 1).Value = 'Commercial and public services'
$ws.Cells.Item(35, 2).Value = 2010
$ws.Cells.Item(35, 3).Value = 0.00465182109398443
$ws.Cells.Item(35, 4).Value = 0.248780292784875
$ws.Cells.Item(35, 5).Value = 0.163919367208958
$ws.Cells.Item(35, 6).Value = 0.025323313743574
$ws.Cells.Item(35, 7).Value = 0.557325205168608
$ws.Cells.Item(36, 1).Value = 'Commercial and public services'
$ws.Cells.Item(36, 2).Value = 2011
$ws.Cells.Item(36, 3).Value = 0.00443937837618143
$ws.Cells.Item(36, 4).Value = 0.247661613325155
$ws.Cells.Item(36, 5).Value = 0.170166010212602
$ws.Cells.Item(36, 6).Value = 0.0269752739741212
$ws.Cells.Item(36, 7).Value = 0.55075772411194
$ws.Cells.Item(37, 1).Value = 'Commercial and public services'
$ws.Cells.Item(37, 2).Value = 2012
$ws.Cells.Item(37, 3).Value = 0.00370231464474871
$ws.Cells.Item(37, 4).Value = 0.241225779086881
$ws.Cells.Item(37, 5).Value = 0.168257442302387
$ws.Cells.Item(37, 6).Value = 0.028215585444752
$ws.Cells.Item(37, 7).Value = 0.558598878521231
$ws.Cells.Item(38, 1).Value = 'Commercial and public services'
$ws.Cells.Item(38, 2).Value = 2013
$ws.Cells.Item(38, 3).Value = 0.00377392792263585
$ws.Cells.Item(38, 4).Value = 0.267811848087296
$ws.Cells.Item(38, 5).Value = 0.166527881251453
$ws.Cells.Item(38, 6).Value = 0.0293294768802379
$ws.Cells.Item(38, 7).Value = 0.532556865858377
$ws.Cells.Item(39, 1).Value = 'Commercial and public services'
$ws.Cells.Item(39, 2).Value = 2014
$ws.Cells.Item(39, 3).Value = 0.0036354345163539
$ws.Cells.Item(39, 4).Value = 0.255314983044822
$ws.Cells.Item(39, 5).Value = 0.167067980866524
$ws.Cells.Item(39, 6).Value = 0.0296629979050195
$ws.Cells.Item(39, 7).Value = 0.54431860366728
$ws.Cells.Item(40, 1).Value = 'Commercial and public services'
$ws.Cells.Item(40, 2).Value = 2015
$ws.Cells.Item(40, 3).Value = 0.00384221225809996
$ws.Cells.Item(40, 4).Value = 0.235266941064898
$ws.Cells.Item(40, 5).Value = 0.186677432377876
$ws.Cells.Item(40, 6).Value = 0.030863828882524
$ws.Cells.Item(40, 7).Value = 0.543349585416602
$ws.Cells.Item(41, 1).Value = 'Commercial and public services'
$ws.Cells.Item(41, 2).Value = 2016
$ws.Cells.Item(41, 3).Value = 0.00344430509299673
$ws.Cells.Item(41, 4).Value = 0.22954751401036
$ws.Cells.Item(41, 5).Value = 0.172581337311204
$ws.Cells.Item(41, 6).Value = 0.0318521241103312
$ws.Cells.Item(41, 7).Value = 0.562574719475109
$ws.Cells.Item(42, 1).Value = 'Commercial and public services'
$ws.Cells.Item(42, 2).Value = 2017
$ws.Cells.Item(42, 3).Value = 0.00352644182596447
$ws.Cells.Item(42, 4).Value = 0.216322072075319
$ws.Cells.Item(42, 5).Value = 0.186930731705442
$ws.Cells.Item(42, 6).Value = 0.0336621891716418
$ws.Cells.Item(42, 7).Value = 0.559558565221633
$ws.Cells.Item(43, 1).Value = 'Commercial and public services'
$ws.Cells.Item(43, 2).Value = 2018
$ws.Cells.Item(43, 3).Value = 0.00227920610692883
$ws.Cells.Item(43, 4).Value = 0.240628984112254
$ws.Cells.Item(43, 5).Value = 0.176965719300653
$ws.Cells.Item(43, 6).Value = 0.0340575770647568
$ws.Cells.Item(43, 7).Value = 0.546068513415407
$ws.Cells.Item(44, 1).Value = 'Commercial and public services'
$ws.Cells.Item(44, 2).Value = 2019
$ws.Cells.Item(44, 3).Value = 0.00254169586888653
$ws.Cells.Item(44, 4).Value = 0.215081088233997
$ws.Cells.Item(44, 5).Value = 0.188233506667908
$ws.Cells.Item(44, 6).Value = 0.0355763905566147
$ws.Cells.Item(44, 7).Value = 0.558567318672594
$ws.Cells.Item(45, 1).Value = 'Commercial and public services'
$ws.Cells.Item(45, 2).Value = 2020
$ws.Cells.Item(45, 3).Value = 0.00269743384889377
$ws.Cells.Item(45, 4).Value = 0.235100657037878
$ws.Cells.Item(45, 5).Value = 0.155260048019108
$ws.Cells.Item(45, 6).Value = 0.0375768243307909
$ws.Cells.Item(45, 7).Value = 0.569365036763329

$ws = $wb.Worksheets.Item("IND_BF")
$ws.Cells.Clear()
$ws.Cells.Item(1, 1).Value = 'Year'
$ws.Cells.Item(1, 2).Value = 'intensity'
$ws.Cells.Item(2, 1).Value = 2010
$ws.Cells.Item(2, 2).Value = 0.154347028848067
$ws.Cells.Item(3, 1).Value = 2011
$ws.Cells.Item(3, 2).Value = 0.144116657652733
$ws.Cells.Item(4, 1).Value = 2012
$ws.Cells.Item(4, 2).Value = 0.145612538466151
$ws.Cells.Item(5, 1).Value = 2013
$ws.Cells.Item(5, 2).Value = 0.151571753065132
$ws.Cells.Item(6, 1).Value = 2014
$ws.Cells.Item(6, 2).Value = 0.14932172493144
$ws.Cells.Item(7, 1).Value = 2015
$ws.Cells.Item(7, 2).Value = 0.140708085157738
$ws.Cells.Item(8, 1).Value = 2016
$ws.Cells.Item(8, 2).Value = 0.139336403700608
$ws.Cells.Item(9, 1).Value = 2017
$ws.Cells.Item(9, 2).Value = 0.134461525118918
$ws.Cells.Item(10, 1).Value = 2018
$ws.Cells.Item(10, 2).Value = 0.132622725426858
$ws.Cells.Item(11, 1).Value = 2019
$ws.Cells.Item(11, 2).Value = 0.127302557927408
$ws.Cells.Item(12, 1).Value = 2020
$ws.Cells.Item(12, 2).Value = 0.102806674705956

Write-Host "done"